$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}

Replace-Text "2025-10-18 Saturday" "2025-10-19 Sunday"

Replace-Text "497×2=994" "773×4=3092"
Replace-Text "254×4=1016" "698×4=2792"
Replace-Text "284×6=1704" "747×5=3735"
Replace-Text "245×7=1715" "649×5=3245"
Replace-Text "994×5=4970" "558×5=2790"

Replace-Text "578×3=1734" "889×3=2667"
Replace-Text "754×8=6032" "585×4=2340"
Replace-Text "338×9=3042" "514×2=1028"
Replace-Text "441×8=3528" "641×7=4487"
Replace-Text "742×5=3710" "419×7=2933"

Replace-Text "684×8=5472" "978×4=3912"
Replace-Text "193×9=1737" "855×4=3420"
Replace-Text "435×9=3915" "990×6=5940"
Replace-Text "609×8=4872" "766×6=4596"
Replace-Text "395×7=2765" "545×4=2180"

Replace-Text "646×5=3230" "354×5=1770"
Replace-Text "106×2=212" "329×9=2961"
Replace-Text "636×2=1272" "130×2=260"
Replace-Text "862×6=5172" "110×3=330"
Replace-Text "911×2=1822" "760×5=3800"

Replace-Text "686×3=2058" "107×6=642"
Replace-Text "933×4=3732" "751×4=3004"
Replace-Text "549×3=1647" "136×7=952"
Replace-Text "817×4=3268" "514×3=1542"
Replace-Text "345×9=3105" "543×8=4344"
